$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write column B (Classifier) for all rows first to control shared-string order
$ws.Range("B2").Value = "StackingCV"
$ws.Range("B3").Value = "Stacking (SGD)"
$ws.Range("B4").Value = "Random Forest"
$ws.Range("B5").Value = "Voting"
$ws.Range("B6").Value = "XGB"
$ws.Range("B7").Value = "Stacking (SVC)"
$ws.Range("B8").Value = "Vecstack"
$ws.Range("B9").Value = "Decision Tree"
$ws.Range("B10").Value = "Stacking (Logistic)"
$ws.Range("B11").Value = "SVM (SVC)"
$ws.Range("B12").Value = "Logistic Regression"
$ws.Range("B13").Value = "Gaussian Naive-Bayes"
$ws.Range("B14").Value = "Bernoulli Naive-Bayes"

# Write column C (Features) for all rows
$ws.Range("C2").Value = "Reduced"
$ws.Range("C3").Value = "Reduced"
$ws.Range("C4").Value = "Reduced"
$ws.Range("C5").Value = "Reduced"
$ws.Range("C6").Value = "Reduced"
$ws.Range("C7").Value = "Reduced"
$ws.Range("C8").Value = "Reduced"
$ws.Range("C9").Value = "Reduced"
$ws.Range("C10").Value = "Reduced"
$ws.Range("C11").Value = "Full"
$ws.Range("C12").Value = "Full"
$ws.Range("C13").Value = "Full"
$ws.Range("C14").Value = "Full"

# Write column D (Parameters) for all rows
$ws.Range("D2").Value = "Best"
$ws.Range("D3").Value = "Best"
$ws.Range("D4").Value = "Best"
$ws.Range("D5").Value = "Best"
$ws.Range("D6").Value = "Default"
$ws.Range("D7").Value = "Best"
$ws.Range("D8").Value = "Best"
$ws.Range("D9").Value = "Best"
$ws.Range("D10").Value = "Best"
$ws.Range("D11").Value = "Default"
$ws.Range("D12").Value = "Default"
$ws.Range("D13").Value = "Default"
$ws.Range("D14").Value = "Default"

# Write remaining numeric columns
$ws.Range("A2").Value = 9
$ws.Range("E2").Value = 0.943471598993105
$ws.Range("F2").Value = 0.9334283834174161
$ws.Range("G2").Value = 0.9523329976502182
$ws.Range("H2").Value = 0.9808645226303815
$ws.Range("I2").Value = 8511
$ws.Range("J2").Value = 8730
$ws.Range("K2").Value = 426
$ws.Range("L2").Value = 607
$ws.Range("A3").Value = 7
$ws.Range("E3").Value = 0.9413921418408668
$ws.Range("F3").Value = 0.928054397894275
$ws.Range("G3").Value = 0.9532499718373324
$ws.Range("H3").Value = 0.9801420883286375
$ws.Range("I3").Value = 8462
$ws.Range("J3").Value = 8741
$ws.Range("K3").Value = 415
$ws.Range("L3").Value = 656
$ws.Range("A4").Value = 4
$ws.Range("E4").Value = 0.9391485170187152
$ws.Range("F4").Value = 0.9290414564597499
$ws.Range("G4").Value = 0.9479632945389436
$ws.Range("H4").Value = 0.9800656968185005
$ws.Range("I4").Value = 8471
$ws.Range("J4").Value = 8691
$ws.Range("K4").Value = 465
$ws.Range("L4").Value = 647
$ws.Range("A5").Value = 11
$ws.Range("E5").Value = 0.9387654591222502
$ws.Range("F5").Value = 0.9264093002851502
$ws.Range("G5").Value = 0.9496346261944912
$ws.Range("H5").Value = 0.9781159135727476
$ws.Range("I5").Value = 8447
$ws.Range("J5").Value = 8708
$ws.Range("K5").Value = 448
$ws.Range("L5").Value = 671
$ws.Range("A6").Value = 12
$ws.Range("E6").Value = 0.9404071358213856
$ws.Range("F6").Value = 0.911274402281202
$ws.Range("G6").Value = 0.967400162999185
$ws.Range("H6").Value = 0.9751876601915893
$ws.Range("I6").Value = 8309
$ws.Range("J6").Value = 8876
$ws.Range("K6").Value = 280
$ws.Range("L6").Value = 809
$ws.Range("A7").Value = 8
$ws.Range("E7").Value = 0.9466455072781
$ws.Range("F7").Value = 0.9224610660232507
$ws.Range("G7").Value = 0.9691208664592695
$ws.Range("H7").Value = 0.9691942715818264
$ws.Range("I7").Value = 8411
$ws.Range("J7").Value = 8888
$ws.Range("K7").Value = 268
$ws.Range("L7").Value = 707
$ws.Range("A8").Value = 10
$ws.Range("E8").Value = 0.9383824012257853
$ws.Range("F8").Value = 0.9354025005483658
$ws.Range("G8").Value = 0.9407677035076109
$ws.Range("H8").Value = 0.9383762175087831
$ws.Range("I8").Value = 8529
$ws.Range("J8").Value = 8619
$ws.Range("K8").Value = 537
$ws.Range("L8").Value = 589
$ws.Range("A9").Value = 1
$ws.Range("E9").Value = 0.8984349348801576
$ws.Range("F9").Value = 0.8811142794472472
$ws.Range("G9").Value = 0.9123325005677947
$ws.Range("H9").Value = 0.9180263876339639
$ws.Range("I9").Value = 8034
$ws.Range("J9").Value = 8384
$ws.Range("K9").Value = 772
$ws.Range("L9").Value = 1084
$ws.Range("A10").Value = 6
$ws.Range("E10").Value = 0.8984349348801576
$ws.Range("F10").Value = 0.8811142794472472
$ws.Range("G10").Value = 0.9123325005677947
$ws.Range("H10").Value = 0.9179473069989308
$ws.Range("I10").Value = 8034
$ws.Range("J10").Value = 8384
$ws.Range("K10").Value = 772
$ws.Range("L10").Value = 1084
$ws.Range("A11").Value = 5
$ws.Range("E11").Value = 0.7633587786259542
$ws.Range("F11").Value = 0.6872294372294372
$ws.Range("G11").Value = 0.8141025641025641
$ws.Range("H11").Value = 0.8596272774844202
$ws.Range("I11").Value = 635
$ws.Range("J11").Value = 765
$ws.Range("K11").Value = 145
$ws.Range("L11").Value = 289
$ws.Range("A12").Value = 3
$ws.Range("E12").Value = 0.7284119514063697
$ws.Range("F12").Value = 0.657490677780215
$ws.Range("G12").Value = 0.7651563497128271
$ws.Range("H12").Value = 0.7784839295979675
$ws.Range("I12").Value = 5995
$ws.Range("J12").Value = 7316
$ws.Range("K12").Value = 1840
$ws.Range("L12").Value = 3123
$ws.Range("A13").Value = 2
$ws.Range("E13").Value = 0.6958520302068513
$ws.Range("F13").Value = 0.5834612853695986
$ws.Range("G13").Value = 0.751412429378531
$ws.Range("H13").Value = 0.7514114551785527
$ws.Range("I13").Value = 5320
$ws.Range("J13").Value = 7396
$ws.Range("K13").Value = 1760
$ws.Range("L13").Value = 3798
$ws.Range("A14").Value = 0
$ws.Range("E14").Value = 0.6720477180693882
$ws.Range("F14").Value = 0.6668128975652555
$ws.Range("G14").Value = 0.6668128975652555
$ws.Range("H14").Value = 0.731317942627083
$ws.Range("I14").Value = 6080
$ws.Range("J14").Value = 6201
$ws.Range("K14").Value = 3038
$ws.Range("L14").Value = 3038
